# Update reference-level column headers (row 1, columns B:I) from the old
# "min / SD_nedre / SD_D / D_M / M_G / G_SG / SG_ovre / max" class-boundary
# labels to the new "pess / X0 / X20 / X40 / X60 / X80 / X100 / opt" labels
# used by the updated demo / indicator reference levels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "pess"
$ws.Range("C1").Value = "X0"
$ws.Range("D1").Value = "X20"
$ws.Range("E1").Value = "X40"
$ws.Range("F1").Value = "X60"
$ws.Range("G1").Value = "X80"
$ws.Range("H1").Value = "X100"
$ws.Range("I1").Value = "opt"

# Move the active selection from A2 to A3, matching the saved cursor
# position recorded in the updated workbook.
$ws.Range("A3").Select()
